$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "28×16=448" "28×59=1652"
Replace-Text "75×82=6150" "84×17=1428"
Replace-Text "97×91=8827" "78×47=3666"
Replace-Text "29×37=1073" "33×65=2145"
Replace-Text "52×70=3640" "41×76=3116"
Replace-Text "71×13=923" "21×97=2037"
Replace-Text "73×77=5621" "61×98=5978"
Replace-Text "31×99=3069" "65×30=1950"
Replace-Text "14×64=896" "91×91=8281"
Replace-Text "77×33=2541" "93×32=2976"
Replace-Text "43×80=3440" "40×21=840"
Replace-Text "56×47=2632" "95×49=4655"
Replace-Text "69×29=2001" "86×31=2666"
Replace-Text "98×86=8428" "11×65=715"
Replace-Text "85×14=1190" "66×37=2442"
Replace-Text "27×76=2052" "84×49=4116"
Replace-Text "67×91=6097" "76×78=5928"
Replace-Text "28×75=2100" "76×96=7296"
Replace-Text "64×49=3136" "26×60=1560"
Replace-Text "84×65=5460" "55×15=825"
Replace-Text "98×95=9310" "94×51=4794"
Replace-Text "41×26=1066" "72×23=1656"
Replace-Text "33×78=2574" "25×15=375"
Replace-Text "30×53=1590" "60×47=2820"
Replace-Text "62×11=682" "15×64=960"

Write-Output "Done"
